# Vignola.xlsx update: "aggiornato a 2/3, aggiornati i report"
#
# A new day (2021-02-08, serial 44235) was inserted into the historical
# series at sheet row 93, shifting every following row down by one, and
# two more days (2021-03-01 / 44256 and 2021-03-02 / 44257) were appended
# at the bottom. Columns C ("somma mobile 7gg.") and D ("somma mobile 7gg.
# per 100mila abitanti") are a centered 7-day rolling window over column B
# and were recomputed for every row whose window membership changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 93 (pushes old rows 93-113 down to 94-114).
$ws.Rows.Item(93).Insert()

# Make the new row's date cell (A93) match the date-column style/number
# format used throughout column A (copy the format from the row above it,
# format-only so no value is touched).
$ws.Cells.Item(92, 1).Copy()
$ws.Cells.Item(93, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. Refresh the full A:D block from row 90 through the (about to be
#    extended) end of the table with the final, correct values - this
#    covers the inserted row, the new rolling-sum values around it, and
#    the two brand-new trailing rows.
$data = @(
    @(90,  44232, 19, 102, 399.2328466867588),
    @(91,  44233, 19,  94, 367.9204665544639),
    @(92,  44234, 18,  94, 367.9204665544639),
    @(93,  44235, 16, 103, 403.1468942032957),
    @(94,  44236,  3, 103, 403.1468942032957),
    @(95,  44237,  5,  87, 340.522133938706),
    @(96,  44238, 23,  81, 324.8659438725587),
    @(97,  44239, 19,  81, 317.0378488394849),
    @(98,  44240,  3,  82, 320.9518963560217),
    @(99,  44241, 12,  78, 305.2957062898744),
    @(100, 44242, 16,  71, 277.8973736741164),
    @(101, 44243,  4,  71, 277.8973736741164),
    @(102, 44244,  1,  92, 360.0923715213902),
    @(103, 44245, 16,  91, 356.1783240048534),
    @(104, 44246, 19,  96, 375.7485615875377),
    @(105, 44247, 24, 107, 418.803084269443),
    @(106, 44248, 11, 109, 426.6311793025167),
    @(107, 44249, 21, 118, 461.8576069513484),
    @(108, 44250, 15, 126, 493.1699870836432),
    @(109, 44251,  3, 127, 497.0840346001801),
    @(110, 44252, 25, 139, 544.0526047986223),
    @(111, 44253, 27, 158, 618.4195076128225),
    @(112, 44254, 25, 152, 594.9352225136013)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}

# 3. Rows 113 and 114 (old rows 112/113, now shifted) keep their B value
#    but have no 7-day window yet (C/D stay blank), matching the source.
$ws.Cells.Item(113, 1).Value = 44255
$ws.Cells.Item(113, 2).Value = 23

$ws.Cells.Item(114, 1).Value = 44256
$ws.Cells.Item(114, 2).Value = 40

# 4. Append the new trailing row 115 (2021-03-02, serial 44257) - also
#    with no completed 7-day window yet.
$ws.Cells.Item(115, 1).Value = 44257
$ws.Cells.Item(115, 2).Value = 9

$ws.Cells.Item(114, 1).Copy()
$ws.Cells.Item(115, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
